$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 115192
$ws.Range("J63").Value = 115192
$ws.Range("L63").Value = 115192
$ws.Range("N63").Value = -116440

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 115192
$ws.Range("J66").Value = 115192
$ws.Range("L66").Value = 345576
$ws.Range("N66").Value = -351816

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 500004000
$ws.Range("I86").Value = 333338660
$ws.Range("K86").Value = 333338660
$ws.Range("M86").Value = -333337537

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2162.842
$ws.Range("I88").Value = 2892.2856
$ws.Range("K88").Value = 2892.2856
$ws.Range("M88").Value = -2486.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 500004000
$ws.Range("I89").Value = 333338660
$ws.Range("K89").Value = 1666693300
$ws.Range("M89").Value = -1666687684

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2162.842
$ws.Range("I91").Value = 2892.2856
$ws.Range("K91").Value = 2892.2856
$ws.Range("M91").Value = -1488.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1417.3334
$ws.Range("I98").Value = 1458.2307
$ws.Range("K98").Value = 1458.2307
$ws.Range("M98").Value = 39.76929999999993

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3337.2222
$ws.Range("I100").Value = 3279.5
$ws.Range("K100").Value = 3279.5
$ws.Range("M100").Value = -2738.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 440633
$ws.Range("I112").Value = 962.25
$ws.Range("K112").Value = 2886.75
$ws.Range("M112").Value = -1778.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1417.3334
$ws.Range("I122").Value = 1458.2307
$ws.Range("K122").Value = 4374.6921
$ws.Range("M122").Value = -1924.6921

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 72166.5
$ws.Range("I137").Value = 104625
$ws.Range("K137").Value = 313875
$ws.Range("M137").Value = -311325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4363
$ws.Range("I141").Value = 2948.8
$ws.Range("K141").Value = 8846.400000000001
$ws.Range("M141").Value = -3666.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7717395
$ws.Range("I32").Value = 3269006.5
$ws.Range("K32").Value = 3269006.5
$ws.Range("M32").Value = -3268719.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2996.6667
$ws.Range("I88").Value = 495
$ws.Range("J88").Value = 8000
$ws.Range("K88").Value = 495
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = -89
$ws.Range("N88").Value = -8812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2996.6667
$ws.Range("I91").Value = 495
$ws.Range("J91").Value = 8000
$ws.Range("K91").Value = 495
$ws.Range("L91").Value = 8000
$ws.Range("M91").Value = 909
$ws.Range("N91").Value = -10808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2672.2888
$ws.Range("I132").Value = 2210.2058
$ws.Range("J132").Value = 4100.5454
$ws.Range("K132").Value = 6630.617400000001
$ws.Range("L132").Value = 12301.6362
$ws.Range("M132").Value = -4100.617400000001
$ws.Range("N132").Value = -17361.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2685.6
$ws.Range("I86").Value = 1981.2
$ws.Range("K86").Value = 1981.2
$ws.Range("M86").Value = -858.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2685.6
$ws.Range("I89").Value = 1981.2
$ws.Range("K89").Value = 9906
$ws.Range("M89").Value = -4290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 69596.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 69596.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 818
$ws.Range("I70").Value = 818
$ws.Range("K70").Value = 2454
$ws.Range("M70").Value = -2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 818
$ws.Range("I73").Value = 818
$ws.Range("K73").Value = 2454
$ws.Range("M73").Value = -1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 922.5
$ws.Range("I129").Value = 563.3333
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 1689.9999
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 3310.0001
$ws.Range("N129").Value = -16000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2197.6365
$ws.Range("I131").Value = 1082
$ws.Range("K131").Value = 3246
$ws.Range("M131").Value = 1794

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3598.75
$ws.Range("I122").Value = 3965
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 11895
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -9445
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1698.6428
$ws.Range("I132").Value = 1698.6428
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5095.928400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2565.928400000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 29000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 29000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 29000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -29876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2970.3845
$ws.Range("I46").Value = 604.375
$ws.Range("J46").Value = 4021.9443
$ws.Range("K46").Value = 604.375
$ws.Range("L46").Value = 4021.9443
$ws.Range("M46").Value = -416.375
$ws.Range("N46").Value = -4397.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 388.64285
$ws.Range("I55").Value = 318
$ws.Range("K55").Value = 318
$ws.Range("M55").Value = -145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2035.5714
$ws.Range("I68").Value = 2043.5
$ws.Range("K68").Value = 2043.5
$ws.Range("M68").Value = -1294.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2035.5714
$ws.Range("I71").Value = 2043.5
$ws.Range("K71").Value = 10217.5
$ws.Range("M71").Value = -6473.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4352.5
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 366333
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 366333
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 366333
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -376693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5245.647
$ws.Range("I62").Value = 5073.7144
$ws.Range("J62").Value = 5366
$ws.Range("K62").Value = 5073.7144
$ws.Range("L62").Value = 5366
$ws.Range("M62").Value = -4449.7144
$ws.Range("N62").Value = -6614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5245.647
$ws.Range("I65").Value = 5073.7144
$ws.Range("J65").Value = 5366
$ws.Range("K65").Value = 25368.572
$ws.Range("L65").Value = 26830
$ws.Range("M65").Value = -22248.572
$ws.Range("N65").Value = -33070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6420.3335
$ws.Range("I122").Value = 4248.5
$ws.Range("K122").Value = 12745.5
$ws.Range("M122").Value = -10295.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3374.1562
$ws.Range("I126").Value = 2211.8333
$ws.Range("K126").Value = 6635.499899999999
$ws.Range("M126").Value = -4165.499899999999
